$wb = $excel.ActiveWorkbook

### ------------------------------------------------------------------
### 1) "总计" sheet: the quarter rows (columns B/C/D) all shift down by
###    one, a brand-new row appears at the bottom carrying what used to
###    be the last row's figures, and the freshly vacated row 2 gets the
###    new 2022-Q4 numbers. Column A is a fixed running index (0,1,2,…)
###    tied to row position, so it is NOT shifted - only a new "7" is
###    appended for the new bottom row.
### ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Grow the table by one row: row 9 becomes a new index row, matching
# the look (bold/centered/bordered) of the rest of column A.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$summary.Range("A9").Value2 = 7

# Shift B/C/D down one row at a time, bottom-up so sources aren't
# clobbered before they are read.
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Range("B$dst").Value2 = $summary.Range("B$r").Value2
    $summary.Range("C$dst").Value2 = $summary.Range("C$r").Value2
    $summary.Range("D$dst").Value2 = $summary.Range("D$r").Value2
}

# Write the new 2022-Q4 figures into the now-empty row 2.
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 13
$summary.Range("D2").Value2 = 1.63

### ------------------------------------------------------------------
### 2) Insert a brand-new worksheet named "2022-Q4" right after "总计"
###    (i.e. right before what is currently "2022-Q3"), holding the
###    fund-holding breakdown for the new quarter.
### ------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value2 = $headers[$i]
}

# data rows: index, code, name, size, stockPosition, positionShare, marketValue, rank
$rows = @(
    @(0,  "000566", "华泰柏瑞创新升级混合A",             "18.24", "89.96", "3.08", "0.5618", 3),
    @(1,  "007968", "华泰柏瑞研究精选混合A",             "8.64",  "91.41", "3.07", "0.2652", 6),
    @(2,  "009636", "华泰柏瑞景气优选混合A",             "8.12",  "91.93", "3.07", "0.2493", 3),
    @(3,  "000967", "华泰柏瑞创新动力灵活配置混合",       "5.51",  "89.56", "3.08", "0.1697", 3),
    @(4,  "013431", "华泰柏瑞景气汇选三年持有期混合A",    "4.96",  "91.63", "3.08", "0.1528", 3),
    @(5,  "008373", "华泰柏瑞景气回报一年持有期混合A",    "2.45",  "89.95", "3.08", "0.0755", 3),
    @(6,  "013847", "华泰柏瑞匠心汇选混合A",             "1.88",  "92.25", "3.10", "0.0583", 3),
    @(7,  "013432", "华泰柏瑞景气汇选三年持有期混合C",    "1.31",  "91.63", "3.08", "0.0403", 3),
    @(8,  "010291", "华泰柏瑞研究精选混合C",             "0.93",  "91.41", "3.07", "0.0286", 6),
    @(9,  "010028", "华泰柏瑞创新升级混合C",             "0.58",  "89.96", "3.08", "0.0179", 3),
    @(10, "013848", "华泰柏瑞匠心汇选混合C",             "0.17",  "92.25", "3.10", "0.0053", 3),
    @(11, "008374", "华泰柏瑞景气回报一年持有期混合C",    "0.07",  "89.95", "3.08", "0.0022", 3),
    @(12, "011454", "华泰柏瑞景气优选混合C",             "0.02",  "91.93", "3.07", "0.0006", 3)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $xlRow = $r + 2

    $q4.Cells.Item($xlRow, 1).Value2 = $row[0]                 # A: index (number)
    $q4.Cells.Item($xlRow, 2).Value2 = "'" + $row[1]           # B: fund code (keep as text)
    $q4.Cells.Item($xlRow, 3).Value2 = $row[2]                 # C: fund name (text)
    $q4.Cells.Item($xlRow, 4).Value2 = "'" + $row[3]           # D: fund size (text)
    $q4.Cells.Item($xlRow, 5).Value2 = "'" + $row[4]           # E: stock position (text)
    $q4.Cells.Item($xlRow, 6).Value2 = "'" + $row[5]           # F: position share (text)
    $q4.Cells.Item($xlRow, 7).Value2 = "'" + $row[6]           # G: market value held (text)
    $q4.Cells.Item($xlRow, 8).Value2 = $row[7]                 # H: rank (number)

    $q4.Range("B" + $xlRow + ":G" + $xlRow).Style = "Normal"
}

# Header row (B1:H1) and index column (A2:A14) use the bold/centered/
# bordered look (style index 2 in the original file) - copy it across
# from the equivalent cells on the "总计" sheet.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A3").Copy()
$q4.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the original active sheet/selection ("总计"!A1), matching the
# untouched bookViews/sheetView state from the source file.
$summary.Activate()
$summary.Range("A1").Select() | Out-Null
